$wb = $excel.ActiveWorkbook
$wsFB = $wb.Worksheets.Item("FB")
$wsFBrieng = $wb.Worksheets.Item("FB rieng")

# Copy rows 2:7 (B:C data) from FB to FB rieng rows 8:13
$wsFB.Range("B2:C7").Copy()
$wsFBrieng.Range("B8").PasteSpecial()

